$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new D text value or $null, new E text value)
$updates = @{
    2  = @('37.524.25', '  +0.17%  ')
    3  = @('2.023.35', '  +0.87%  ')
    4  = @($null, '  -0.02%  ')
    5  = @('254.20', '  +2.99%  ')
    6  = @('0.619', '  -2.17%  ')
    7  = @($null, '  -0.04%  ')
    8  = @('56.69', '  -8.63%  ')
    9  = @($null, '  -0.93%  ')
    10 = @('0.0782', '  -2.80%  ')
    11 = @($null, '  -2.32%  ')
    12 = @('14.49', '  -3.09%  ')
    13 = @('2.324.00', '  +1.11%  ')
    14 = @('0.814', '  -4.76%  ')
    15 = @('20.99', '  -7.86%  ')
    16 = @('5.33', '  -2.23%  ')
    17 = @('2.032.03', '  +1.14%  ')
    18 = @('37.427.22', '  +0.24%  ')
    19 = @($null, '  -1.45%  ')
    20 = @($null, '  -2.58%  ')
    21 = @($null, '  -0.75%  ')
    22 = @('228.31', '  -1.60%  ')
    23 = @($null, '  +3.89%  ')
    24 = @($null, '  -0.07%  ')
    25 = @('2.32', '  -1.64%  ')
    26 = @('163.40', '  -0.37%  ')
    27 = @('9.02', '  -3.69%  ')
    28 = @('19.75', '  -0.08%  ')
    29 = @($null, '  -9.08%  ')
    30 = @($null, '  +0.25%  ')
    31 = @('0.120', '  -1.24%  ')
    32 = @('0.0666', '  +6.21%  ')
    33 = @($null, '  -4.44%  ')
    34 = @($null, '  -1.09%  ')
    35 = @($null, '  +3.02%  ')
    36 = @($null, '  +0.13%  ')
    37 = @($null, '  +0.63%  ')
    38 = @($null, '  +0.20%  ')
    39 = @($null, '  -3.41%  ')
    40 = @($null, '  +2.94%  ')
    41 = @('0.0965', '  -1.93%  ')
    42 = @($null, '  +1.18%  ')
    43 = @('0.0215', '  -0.42%  ')
    44 = @('1.408.75', '  +1.58%  ')
    45 = @('15.93', '  -4.99%  ')
    46 = @('90.66', '  -0.90%  ')
    47 = @($null, '  -1.08%  ')
    48 = @('7.30', '  +0.05%  ')
    49 = @($null, '  +0.71%  ')
    50 = @('2.01', '  +0.57%  ')
    51 = @('2.215.86', '  +1.17%  ')
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $newD = $pair[0]
    $newE = $pair[1]
    if ($null -ne $newD) {
        # Leading apostrophe forces Excel to treat the value as literal text,
        # so numeric-looking strings (e.g. "254.20") are not coerced into
        # numbers and do not lose formatting (trailing zeros etc.).
        $ws.Range("D$row").Value = "'" + $newD
        $ws.Range("D$row").Style = "Normal"
    }
    $ws.Range("E$row").Value = $newE
    $ws.Range("E$row").Style = "Normal"
}
